$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the mislabeled strain: every "KN99alpha" (col F) should read "TDY451"
[void]$ws.Cells.Replace("KN99alpha", "TDY451")

# Move the selection to F17, matching the resulting saved view state
$ws.Range("F17").Select() | Out-Null

# Column G (genotype) now needs to fit its (wider) content, e.g. "CNAG_00000"
$ws.Columns.Item(7).ColumnWidth = 10.5
